# optimize: skip existed res in multi-process
#
# The sweep table on Sheet1 (backed by table "表1") is extended from a
# 4-run grid (A1:R5) to a 20-run grid (A1:R21): for every mkt_type in
# (CSI500, CSI300) and every "G" value in (1.0, 1.5, 2.0, 2.5, 0.5) there
# are two rows - one per alpha_name. The first 4 existing rows are
# rewritten in place (their mkt_type flips from CSI300 to CSI500, and
# their "D"/"G" columns settle on the new grid values); 16 new rows are
# appended below them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns, in sheet order:
# A=run B=alpha_name C=mkt_type D=beta_kind E=beta_suffix F=beta_args
# G=H0 H=H1 I=B J=E K=D L=G M=S N=N O=wei_tole P=begin_date Q=end_date
# R=opt_verbose

$alphaFRtn = "FRtn5D(0.0,3.0)_zscore_SD(0.0225)"
$alphaFactorApm = "factor_apm_zscore_SD(0.0225)"
$betaArgs = "(['size', 'beta', 'momentum'],)"

$gVals = @("1.0", "1.5", "2.0", "2.5", "0.5")
$mktTypes = @("CSI500", "CSI300")

# "FALSE" is written with a leading apostrophe so it is stored as the
# literal text "FALSE" (shared string) instead of being auto-coerced into
# a real boolean cell by the COM value-assignment heuristics.
$optVerbose = "'FALSE"

# Build the full 20-row grid (mkt_type outer, G middle, alpha_name inner)
# matching the row order already on the sheet.
$allRows = New-Object System.Collections.ArrayList
foreach ($mkt in $mktTypes) {
    foreach ($g in $gVals) {
        foreach ($alpha in @($alphaFRtn, $alphaFactorApm)) {
            $row = @("1", $alpha, $mkt, "Barra", "barra3", $betaArgs, "0.20", "0.02", "0", "0.5", "2", $g, "inf", "inf", "1e-5", "2016-02-01", "2022-03-31", $optVerbose)
            [void]$allRows.Add($row)
        }
    }
}

$lo = $ws.ListObjects.Item(1)

# Grow the table from 4 data rows to 20 data rows.
while ($lo.ListRows.Count -lt $allRows.Count) {
    [void]$lo.ListRows.Add()
}

# Write every data row (this both rewrites the original 4 rows in place
# and fills in the 16 newly-added ones).
for ($i = 0; $i -lt $allRows.Count; $i++) {
    $r = 2 + $i
    $vals = $allRows[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($r, $col).Value = $vals[$c]
    }
}

# Match the saved selection from the edit.
$ws.Range("C12:C21").Select()
